$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1844.2285
$ws.Range("I70").Value = 1033.5
$ws.Range("J70").Value = 2925.2
$ws.Range("K70").Value = 3100.5
$ws.Range("L70").Value = 8775.599999999999
$ws.Range("M70").Value = -2830.5
$ws.Range("N70").Value = -9315.599999999999
$ws.Range("H73").Value = 1844.2285
$ws.Range("I73").Value = 1033.5
$ws.Range("J73").Value = 2925.2
$ws.Range("K73").Value = 3100.5
$ws.Range("L73").Value = 8775.599999999999
$ws.Range("M73").Value = -2164.5
$ws.Range("N73").Value = -10647.6
$ws.Range("H112").Value = 1103.9333
$ws.Range("J112").Value = 1185.3077
$ws.Range("L112").Value = 3555.9231
$ws.Range("N112").Value = -5771.9231
$ws.Range("H113").Value = 8308.166999999999
$ws.Range("I113").Value = 3012.4
$ws.Range("J113").Value = 12090.857
$ws.Range("K113").Value = 3012.4
$ws.Range("L113").Value = 12090.857
$ws.Range("M113").Value = 241.5999999999999
$ws.Range("N113").Value = -18598.857
$ws.Range("H137").Value = 1999.05
$ws.Range("I137").Value = 1317.1875
$ws.Range("J137").Value = 4726.5
$ws.Range("K137").Value = 3951.5625
$ws.Range("L137").Value = 14179.5
$ws.Range("M137").Value = -1401.5625
$ws.Range("N137").Value = -19279.5
$ws.Range("H138").Value = 1704.6056
$ws.Range("I138").Value = 1167.0927
$ws.Range("J138").Value = 3412
$ws.Range("K138").Value = 3501.2781
$ws.Range("L138").Value = 10236
$ws.Range("M138").Value = 1638.7219
$ws.Range("N138").Value = -20516
$ws.Range("H141").Value = 6606.1035
$ws.Range("I141").Value = 974.5122
$ws.Range("J141").Value = 20188.176
$ws.Range("K141").Value = 2923.5366
$ws.Range("L141").Value = 60564.528
$ws.Range("M141").Value = 2256.4634
$ws.Range("N141").Value = -70924.52799999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5492.7793
$ws.Range("I32").Value = 3900.1428
$ws.Range("K32").Value = 3900.1428
$ws.Range("M32").Value = -3613.1428
$ws.Range("H61").Value = 1250
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -788
$ws.Range("N61").Value = -1924
$ws.Range("H74").Value = 1272.931
$ws.Range("I74").Value = 1346.28
$ws.Range("J74").Value = 814.5
$ws.Range("K74").Value = 1346.28
$ws.Range("L74").Value = 814.5
$ws.Range("M74").Value = -472.28
$ws.Range("N74").Value = -2562.5
$ws.Range("H77").Value = 1272.931
$ws.Range("I77").Value = 1346.28
$ws.Range("J77").Value = 814.5
$ws.Range("K77").Value = 6731.4
$ws.Range("L77").Value = 4072.5
$ws.Range("M77").Value = -2363.4
$ws.Range("N77").Value = -12808.5
$ws.Range("H132").Value = 2241.8948
$ws.Range("I132").Value = 1334.3846
$ws.Range("J132").Value = 4208.1665
$ws.Range("K132").Value = 4003.1538
$ws.Range("L132").Value = 12624.4995
$ws.Range("M132").Value = -1473.1538
$ws.Range("N132").Value = -17684.4995
$ws.Range("H136").Value = 1250
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -9600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 12000
$ws.Range("J32").Value = 12000
$ws.Range("L32").Value = 12000
$ws.Range("N32").Value = -12768
$ws.Range("H75").Value = 12655.444
$ws.Range("I75").Value = 8299.666999999999
$ws.Range("K75").Value = 8299.666999999999
$ws.Range("M75").Value = -7363.666999999999
$ws.Range("H78").Value = 12655.444
$ws.Range("I78").Value = 8299.666999999999
$ws.Range("K78").Value = 24899.001
$ws.Range("M78").Value = -20219.001
$ws.Range("H134").Value = 1691.1724
$ws.Range("I134").Value = 1466.0952
$ws.Range("J134").Value = 2282
$ws.Range("K134").Value = 4398.2856
$ws.Range("L134").Value = 6846
$ws.Range("M134").Value = -1863.2856
$ws.Range("N134").Value = -11916

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2497.1785
$ws.Range("I31").Value = 1644.4688
$ws.Range("J31").Value = 3634.125
$ws.Range("K31").Value = 1644.4688
$ws.Range("L31").Value = 3634.125
$ws.Range("M31").Value = -1349.4688
$ws.Range("N31").Value = -4224.125
$ws.Range("H34").Value = 2497.1785
$ws.Range("I34").Value = 1644.4688
$ws.Range("J34").Value = 3634.125
$ws.Range("K34").Value = 1644.4688
$ws.Range("L34").Value = 3634.125
$ws.Range("M34").Value = -1442.4688
$ws.Range("N34").Value = -4038.125
$ws.Range("H94").Value = 1901.2858
$ws.Range("I94").Value = 2642
$ws.Range("J94").Value = 1489.7778
$ws.Range("K94").Value = 2642
$ws.Range("L94").Value = 1489.7778
$ws.Range("M94").Value = -2191
$ws.Range("N94").Value = -2391.7778
$ws.Range("H99").Value = 1913477
$ws.Range("I99").Value = 5335002
$ws.Range("K99").Value = 5335002
$ws.Range("M99").Value = -5333504
$ws.Range("H126").Value = 1913477
$ws.Range("I126").Value = 5335002
$ws.Range("K126").Value = 16005006
$ws.Range("M126").Value = -16002536
$ws.Range("H132").Value = 2117.054
$ws.Range("J132").Value = 3454.842
$ws.Range("L132").Value = 10364.526
$ws.Range("N132").Value = -15424.526
$ws.Range("H134").Value = 1257.6041
$ws.Range("I134").Value = 1294.659
$ws.Range("K134").Value = 3883.977
$ws.Range("M134").Value = -1348.977
$ws.Range("H141").Value = 33519.23
$ws.Range("J141").Value = 33519.23
$ws.Range("L141").Value = 33519.23
$ws.Range("N141").Value = -43879.23

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7679.256
$ws.Range("J39").Value = 5517.1953
$ws.Range("L39").Value = 16551.5859
$ws.Range("N39").Value = -17139.5859

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2760.6086
$ws.Range("I132").Value = 2101.125
$ws.Range("J132").Value = 4268
$ws.Range("K132").Value = 6303.375
$ws.Range("L132").Value = 12804
$ws.Range("M132").Value = -3773.375
$ws.Range("N132").Value = -17864

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 145629.14
$ws.Range("I7").Value = 201680.8
$ws.Range("J7").Value = 5500
$ws.Range("K7").Value = 201680.8
$ws.Range("L7").Value = 5500
$ws.Range("M7").Value = -201568.8
$ws.Range("N7").Value = -5724
$ws.Range("H100").Value = 720.63635
$ws.Range("I100").Value = 720.63635
$ws.Range("K100").Value = 720.63635
$ws.Range("M100").Value = -179.63635
$ws.Range("H126").Value = 145629.14
$ws.Range("I126").Value = 201680.8
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 605042.3999999999
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -602572.3999999999
$ws.Range("N126").Value = -21440
$ws.Range("H132").Value = 5454.817
$ws.Range("I132").Value = 5565.341
$ws.Range("J132").Value = 5274.7036
$ws.Range("K132").Value = 16696.023
$ws.Range("L132").Value = 15824.1108
$ws.Range("M132").Value = -14166.023
$ws.Range("N132").Value = -20884.1108
$ws.Range("H136").Value = 14495215
$ws.Range("I136").Value = 2768.0527
$ws.Range("J136").Value = 83334340
$ws.Range("K136").Value = 8304.158100000001
$ws.Range("L136").Value = 250003020
$ws.Range("M136").Value = -5754.158100000001
$ws.Range("N136").Value = -250008120

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 2600
$ws.Range("I32").Value = 2600
$ws.Range("K32").Value = 2600
$ws.Range("M32").Value = -2283
$ws.Range("H136").Value = 2690.3809
$ws.Range("I136").Value = 2284.4614
$ws.Range("J136").Value = 3350
$ws.Range("K136").Value = 6853.3842
$ws.Range("L136").Value = 10050
$ws.Range("M136").Value = -4303.3842
$ws.Range("N136").Value = -15150
